# Added a Gantt Chart: close the date gaps in Team -1's schedule (rows 5-16)
# so each task's Start Date follows immediately after the previous task's
# End Date (skipping weekends), and refresh Duration/Status accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates are stored as plain text in this sheet (not real Excel dates), so
# force the Start/End Date columns to Text format before writing the new
# values - otherwise Excel would auto-convert "2025-08-06" into a date
# serial number.
$ws.Range("E5:F16").NumberFormat = "@"

# Row 5 - WARUKIRA: now a 2-day "In Progress" task instead of a 7-day
# "Not Done" one; Start Date is unchanged.
$ws.Range("D5").Value = 2
$ws.Range("F5").Value = "2025-08-05"
$ws.Range("G5").Value = "In Progress"

# Row 6 - CHAKA RELI
$ws.Range("E6").Value = "2025-08-06"
$ws.Range("F6").Value = "2025-08-07"

# Row 7 - JAMHURI
$ws.Range("E7").Value = "2025-08-08"
$ws.Range("F7").Value = "2025-08-15"

# Row 8 - NYAIROKO PASSENGA
$ws.Range("E8").Value = "2025-08-16"
$ws.Range("F8").Value = "2025-08-23"

# Row 9 - GWA KARUMBI BAHATI
$ws.Range("E9").Value = "2025-08-25"
$ws.Range("F9").Value = "2025-08-25"

# Row 10 - KWA CHUI
$ws.Range("E10").Value = "2025-08-26"
$ws.Range("F10").Value = "2025-08-26"

# Row 11 - IGANJO VILLAGE
$ws.Range("E11").Value = "2025-08-27"
$ws.Range("F11").Value = "2025-09-03"

# Row 12 - MAHI MAHORO
$ws.Range("E12").Value = "2025-09-04"
$ws.Range("F12").Value = "2025-09-11"

# Row 13 - CHURIRI
$ws.Range("E13").Value = "2025-09-12"
$ws.Range("F13").Value = "2025-09-12"

# Row 14 - LOWER GITHIMA
$ws.Range("E14").Value = "2025-09-13"
$ws.Range("F14").Value = "2025-09-13"

# Row 15 - MBIRITI
$ws.Range("E15").Value = "2025-09-15"
$ws.Range("F15").Value = "2025-09-15"

# Row 16 - KARIAHU
$ws.Range("E16").Value = "2025-09-16"
$ws.Range("F16").Value = "2025-09-17"
